$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.722.42'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.10%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.227.97'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.10%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '273.87'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +6.28%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '87.70'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +8.85%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.621'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.24%  '

$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('E9').Value = '  -0.45%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '45.19'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.93%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0920'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.54%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.67'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.89%  '

$ws.Range('E13').Value = '  +1.03%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.558.47'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.16%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.94'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.82%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.226.75'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.71%  '

$ws.Range('E17').Value = '  -1.33%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.637.62'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.06%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000104'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.29%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '70.23'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.05%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.96'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.90%  '

$ws.Range('E22').Value = '  -0.43%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '232.39'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.32%  '

$ws.Range('E24').Value = '  -7.66%  '

$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.57'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +14.81%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.82'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.52%  '

$ws.Range('E28').Value = '  +3.39%  '

$ws.Range('E29').Value = '  +4.47%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '39.11'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.66%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '172.83'
$ws.Range('D31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0905'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.01%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.81'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.52%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.36'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.19%  '

$ws.Range('E35').Value = '  +0.57%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.111'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.40%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0354'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.52%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.27'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -6.06%  '

$ws.Range('E39').Value = '  +16.54%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '12.52'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.75%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.17'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.62%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '63.79'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.06%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.207'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.90%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.38'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.38%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.44'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.70%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0986'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.30%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '100.40'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.76%  '

$ws.Range('E48').Value = '  +2.92%  '

$ws.Range('E49').Value = '  +0.27%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.49'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.10%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.428'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -7.33%  '
